# Algerian forest fires dataset update:
# - Normalize "fire"/"not fire" class labels in Region 1 (strip padding variants)
# - Region 2: a few class-label cells are retargeted to the surviving shared-string entries
# - Update the active selection / scroll position on the Region 1 sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Region 1")
$ws2 = $wb.Worksheets.Item("Region 2")

# --- Region 1 (sheet1): normalize every class-label cell in column N ---
$ws1.Range("N1").Value = "class"
$ws1.Range("N2").Value = "not fire"
$ws1.Range("N3").Value = "not fire"
$ws1.Range("N4").Value = "not fire"
$ws1.Range("N5").Value = "not fire"
$ws1.Range("N6").Value = "not fire"
$ws1.Range("N7").Value = "fire"
$ws1.Range("N8").Value = "fire"
$ws1.Range("N9").Value = "fire"
$ws1.Range("N10").Value = "not fire"
$ws1.Range("N11").Value = "not fire"
$ws1.Range("N12").Value = "fire"
$ws1.Range("N13").Value = "fire"
$ws1.Range("N14").Value = "not fire"
$ws1.Range("N15").Value = "not fire"
$ws1.Range("N16").Value = "not fire"
$ws1.Range("N17").Value = "not fire"
$ws1.Range("N18").Value = "not fire"
$ws1.Range("N19").Value = "not fire"
$ws1.Range("N20").Value = "not fire"
$ws1.Range("N21").Value = "not fire"
$ws1.Range("N22").Value = "fire"
$ws1.Range("N23").Value = "not fire"
$ws1.Range("N24").Value = "fire"
$ws1.Range("N25").Value = "fire"
$ws1.Range("N26").Value = "fire"
$ws1.Range("N27").Value = "fire"
$ws1.Range("N28").Value = "fire"
$ws1.Range("N29").Value = "fire"
$ws1.Range("N30").Value = "not fire"
$ws1.Range("N31").Value = "fire"
$ws1.Range("N32").Value = "not fire"
$ws1.Range("N33").Value = "not fire"
$ws1.Range("N34").Value = "not fire"
$ws1.Range("N35").Value = "not fire"
$ws1.Range("N36").Value = "fire"
$ws1.Range("N37").Value = "fire"
$ws1.Range("N38").Value = "not fire"
$ws1.Range("N39").Value = "fire"
$ws1.Range("N40").Value = "not fire"
$ws1.Range("N41").Value = "not fire"
$ws1.Range("N42").Value = "not fire"
$ws1.Range("N43").Value = "not fire"
$ws1.Range("N44").Value = "not fire"
$ws1.Range("N45").Value = "not fire"
$ws1.Range("N46").Value = "not fire"
$ws1.Range("N47").Value = "not fire"
$ws1.Range("N48").Value = "fire"
$ws1.Range("N49").Value = "fire"
$ws1.Range("N50").Value = "fire"
$ws1.Range("N51").Value = "fire"
$ws1.Range("N52").Value = "fire"
$ws1.Range("N53").Value = "not fire"
$ws1.Range("N54").Value = "not fire"
$ws1.Range("N55").Value = "not fire"
$ws1.Range("N56").Value = "fire"
$ws1.Range("N57").Value = "fire"
$ws1.Range("N58").Value = "fire"
$ws1.Range("N59").Value = "fire"
$ws1.Range("N60").Value = "fire"
$ws1.Range("N61").Value = "fire"
$ws1.Range("N62").Value = "fire"
$ws1.Range("N63").Value = "not fire"
$ws1.Range("N64").Value = "not fire"
$ws1.Range("N65").Value = "not fire"
$ws1.Range("N66").Value = "fire"
$ws1.Range("N67").Value = "fire"
$ws1.Range("N68").Value = "fire"
$ws1.Range("N69").Value = "fire"
$ws1.Range("N70").Value = "not fire"
$ws1.Range("N71").Value = "fire"
$ws1.Range("N72").Value = "fire"
$ws1.Range("N73").Value = "fire"
$ws1.Range("N74").Value = "not fire"
$ws1.Range("N75").Value = "fire"
$ws1.Range("N76").Value = "fire"
$ws1.Range("N77").Value = "fire"
$ws1.Range("N78").Value = "fire"
$ws1.Range("N79").Value = "fire"
$ws1.Range("N80").Value = "fire"
$ws1.Range("N81").Value = "fire"
$ws1.Range("N82").Value = "fire"
$ws1.Range("N83").Value = "fire"
$ws1.Range("N84").Value = "fire"
$ws1.Range("N86").Value = "fire"
$ws1.Range("N88").Value = "fire"
$ws1.Range("N90").Value = "fire"
$ws1.Range("N91").Value = "fire"
$ws1.Range("N92").Value = "not fire"
$ws1.Range("N93").Value = "not fire"
$ws1.Range("N94").Value = "not fire"
$ws1.Range("N95").Value = "not fire"
$ws1.Range("N96").Value = "not fire"
$ws1.Range("N97").Value = "not fire"
$ws1.Range("N98").Value = "fire"
$ws1.Range("N99").Value = "not fire"
$ws1.Range("N100").Value = "not fire"
$ws1.Range("N101").Value = "not fire"
$ws1.Range("N102").Value = "not fire"
$ws1.Range("N103").Value = "not fire"
$ws1.Range("N104").Value = "not fire"
$ws1.Range("N105").Value = "not fire"
$ws1.Range("N106").Value = "not fire"
$ws1.Range("N107").Value = "not fire"
$ws1.Range("N108").Value = "not fire"
$ws1.Range("N109").Value = "not fire"
$ws1.Range("N110").Value = "fire"
$ws1.Range("N111").Value = "fire"
$ws1.Range("N112").Value = "fire"
$ws1.Range("N114").Value = "fire"
$ws1.Range("N115").Value = "not fire"
$ws1.Range("N116").Value = "not fire"
$ws1.Range("N117").Value = "not fire"
$ws1.Range("N118").Value = "not fire"
$ws1.Range("N119").Value = "not fire"
$ws1.Range("N120").Value = "fire"
$ws1.Range("N121").Value = "not fire"
$ws1.Range("N122").Value = "not fire"
$ws1.Range("N123").Value = "not fire"

# --- Region 2 (sheet2): a few stray class-label references need retargeting ---
$ws2.Range("N31").Value = "not fire     "
$ws2.Range("J45").Value = "14.6 9"
$ws2.Range("N123").Value = "not fire    "

# --- Update view/selection state on Region 1 to match the saved workbook state ---
$ws1.Activate()
$ws1.Range("N98").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 80
$win.ScrollColumn = 1
